$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.029.03'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").Value = '1.597.21'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '211.85'
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +1.07%  '
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '18.26'
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("D12").Value = '1.818.16'
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").Value = '1.607.73'
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").Value = '  +2.64%  '
$ws.Range("D16").Value = '26.018.03'
$ws.Range("E16").Value = '  +0.55%  '
$ws.Range("D17").Value = '60.87'
$ws.Range("E17").Value = '  +1.31%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("E19").Value = '  -0.15%  '
$ws.Range("D20").Value = '203.61'
$ws.Range("E20").Value = '  +5.02%  '
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").Value = '6.03'
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("E24").Value = '  +13.37%  '
$ws.Range("D25").Value = '143.76'
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -7.58%  '
$ws.Range("E28").Value = '  +0.84%  '
$ws.Range("E29").Value = '  +1.41%  '
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").Value = '3.14'
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("E33").Value = '  -3.80%  '
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").Value = '1.130.90'
$ws.Range("E36").Value = '  +3.17%  '
$ws.Range("E37").Value = '  +7.91%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '0.796'
$ws.Range("E39").Value = '  +2.41%  '
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("E41").Value = '  -1.56%  '
$ws.Range("D42").Value = '0.779'
$ws.Range("E42").Value = '  -2.88%  '
$ws.Range("E43").Value = '  +1.45%  '
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '92.08'
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").Value = '1.51'
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").Value = '54.00'
$ws.Range("E47").Value = '  +1.71%  '
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '0.0₇0950'
$ws.Range("E51").Value = '  -14.41%  '
